$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @("[Bài tập] Xây dựng lớp mô tả hình chữ nhật", "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s16_lap_trinh_hdt/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%5D%20X%C3%A2y%20d%E1%BB%B1ng%20l%E1%BB%9Bp%20m%C3%B4%20t%E1%BA%A3%20h%C3%ACnh%20ch%E1%BB%AF%20nh%E1%BA%ADt.html"),
    @("[Bài tập] Xây dựng lớp mô tả Temperature", "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s16_lap_trinh_hdt/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%5D%20X%C3%A2y%20d%E1%BB%B1ng%20l%E1%BB%9Bp%20m%C3%B4%20t%E1%BA%A3%20Temperature.html"),
    @("[Bài tập] Xây dựng lớp mô tả điện thoại", "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s16_lap_trinh_hdt/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%5D%20X%C3%A2y%20d%E1%BB%B1ng%20l%E1%BB%9Bp%20m%C3%B4%20t%E1%BA%A3%20%C4%91i%E1%BB%87n%20tho%E1%BA%A1i.html")
)

$startRow = 90
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
}
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

$ws.Range("B90:B92").Select()
